$d = $word.ActiveDocument

# --- Portuguese ("Programa") paragraph: split the single run into
#     segments separated by manual line breaks (<w:br/>) at each ';' ---

$d.Content.Find.Execute(
    "\(mínimo 3\); Inovação", $true, $false, $true, $false, $false,
    $true, 1, $false, "(mínimo 3);^l Inovação", 2)

$d.Content.Find.Execute(
    "características;Legislação", $true, $false, $true, $false, $false,
    $true, 1, $false, "características;^lLegislação", 2)

$d.Content.Find.Execute(
    "empresarial;Gerenciamento", $true, $false, $true, $false, $false,
    $true, 1, $false, "empresarial;^lGerenciamento", 2)

$d.Content.Find.Execute(
    "causas;Formulação", $true, $false, $true, $false, $false,
    $true, 1, $false, "causas;^lFormulação", 2)

$d.Content.Find.Execute(
    "etc;Especificação", $true, $false, $true, $false, $false,
    $true, 1, $false, "etc;^lEspecificação", 2)

$d.Content.Find.Execute(
    "Decisão;Elaboração", $true, $false, $true, $false, $false,
    $true, 1, $false, "Decisão;^lElaboração", 2)

# --- English ("Programa" translation, italic) paragraph: same split ---

$d.Content.Find.Execute(
    "\(minimum 3\);Systematic", $true, $false, $true, $false, $false,
    $true, 1, $false, "(minimum 3);^lSystematic", 2)

$d.Content.Find.Execute(
    "characteristics;Legislation", $true, $false, $true, $false, $false,
    $true, 1, $false, "characteristics;^lLegislation", 2)

$d.Content.Find.Execute(
    "business action;Project and Schedule", $true, $false, $true, $false, $false,
    $true, 1, $false, "business action;^lProject and Schedule", 2)

$d.Content.Find.Execute(
    "locate causes;Project Formulation", $true, $false, $true, $false, $false,
    $true, 1, $false, "locate causes;^lProject Formulation", 2)

$d.Content.Find.Execute(
    "\(EAP\), etc;Problem Specification", $true, $false, $true, $false, $false,
    $true, 1, $false, "(EAP), etc;^lProblem Specification", 2)

$d.Content.Find.Execute(
    "decision making;Preparation of reports", $true, $false, $true, $false, $false,
    $true, 1, $false, "decision making;^lPreparation of reports", 2)
